$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.456.48"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.572.33"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'292.18"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.3721"
$ws.Range("D8").Value = "'49.82"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.3413"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'1.151"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'21.29"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "'6.045"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'6.972"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "1.585.66"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'91.33"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "'0.06755"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'6.310"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "'16.35"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "22.455.75"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'2.375"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").Value = "'2.667"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'20.04"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "'149.01"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'5.041"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "'125.81"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "1.756.20"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +10.05%  "
$ws.Range("D33").Value = "'6.224"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "'2.015"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'9.863"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").Value = "'0.08388"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "'0.02492"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "'0.2310"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "'1.346"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'0.06533"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'5.471"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'0.6244"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.06"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'3.811"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'0.5830"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'130.56"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("D49").Value = "'2.080"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("D51").Value = "'0.07330"
$ws.Range("E51").Value = "  -0.02%  "
